# Insert a new data row at row 101 (pushes existing rows 101-162 down to 102-163)
# and populate it with a new "Pepino dulce" price record for Vega Modelo de Temuco.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A101").EntireRow.Insert()

$ws.Range("A101").Value = 10
$ws.Range("B101").Value = "Vega Modelo de Temuco"
$ws.Range("C101").Value = "La Araucanía"
$ws.Range("D101").Value = 44488
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = 100112043
$ws.Range("G101").Value = "Pepino dulce"
$ws.Range("H101").Value = "Cultivar IV Región"
$ws.Range("I101").Value = "Segunda"
$ws.Range("J101").Value = 55
$ws.Range("K101").Value = 16000
$ws.Range("L101").Value = 16000
$ws.Range("M101").Value = 16000
$ws.Range("N101").Value = "$/bandeja 18 kilos"
$ws.Range("O101").Value = "Provincia de Limarí"
$ws.Range("P101").Value = 889
$ws.Range("Q101").Value = 18
$ws.Range("R101").Value = "Hortaliza"
